# issue #5: stock data from json to db
#
# The "股票" (stock) sheet gains a "category" column (constant "normal")
# right after "property_category", plus two trailing columns:
# "source_file" (constant "tmp9edb1") and "index" (mirrors column A's id).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Insert the new "category" column between H (property_category) and
# --- the old I (date) column, shifting date / legislator_name / legislator_id
# --- one column to the right.
$ws.Columns.Item(9).Insert()

$ws.Range("I1").Value = "category"
$ws.Range("I1").Font.Bold = $true

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# --- Append the two trailing columns: source_file, index.
$ws.Range("M1").Value = "source_file"
$ws.Range("M1").Font.Bold = $true
$ws.Range("N1").Value = "index"
$ws.Range("N1").Font.Bold = $true

# Match the header row's border + alignment formatting for the new cells.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmp9edb1"
    $idVal = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 14).Value = $idVal
}

Write-Output "stock sheet updated"
